# Trade #9 closed at 2026-02-17 07:58:18 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets of the live-trading-results workbook to reflect the newly closed
# trade (#9 / row index 8, 0-based) on the MarketMaking strategy.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0.01      # Total P&L $
$summary.Range("B5").Value = 0.02      # Total P&L %
$summary.Range("B6").Value = 9         # Total Trades
$summary.Range("B7").Value = 4         # Winning Trades
$summary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01     # Capital
$status.Range("D4").Value = 9          # Trades
$status.Range("E4").Value = 0.01       # P&L $
$status.Range("F4").Value = 0.01       # P&L %
$status.Range("G4").Value = 44.44      # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the new trade #9 row (spreadsheet row 10) to a trade log
# sheet. The Date column ("2026-02-17") looks like a date to Excel's COM
# layer and would otherwise get silently converted to a date serial number,
# so it is temporarily forced to Text format, written, then the format is
# cleared again so the cell ends up as a plain (General) string cell - same
# as all the other rows above it.
# ---------------------------------------------------------------------------
function Add-Trade9Row($ws) {
    $ws.Range("A10").Value = 9

    $ws.Range("B10").NumberFormat = "@"
    $ws.Range("B10").Value = "2026-02-17"
    $ws.Range("B10").ClearFormats()

    $ws.Range("C10").Value = "07:58:11"
    $ws.Range("D10").Value = "MarketMaking"
    $ws.Range("E10").Value = "DOWN"
    $ws.Range("F10").Value = 0.86
    $ws.Range("G10").Value = 0.88
    $ws.Range("H10").Value = "CLOSED"
    $ws.Range("I10").Value = 2.3256
    $ws.Range("J10").Value = 0.02
    $ws.Range("K10").Value = 100.01
    $ws.Range("L10").Value = 0
    $ws.Range("M10").Value = 0
    $ws.Range("N10").Value = 0.6
    $ws.Range("O10").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P10").Value = "early_exit"
    $ws.Range("Q10").Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade9Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet (same trade log, strategy-specific tab)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade9Row $marketMaking
